$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether it must be forced to text
# (price column D holds numeric-looking strings like "43.362.84" or "6.30" that
# Excel would otherwise silently convert to a Number/Date - keep them Text).
$updates = @(
    @{ Cell = "D2"; Value = "43.362.84"; ForceText = $true }
    @{ Cell = "E2"; Value = "  +0.51%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "2.372.86"; ForceText = $true }
    @{ Cell = "E3"; Value = "  +2.85%  "; ForceText = $false }
    @{ Cell = "E4"; Value = "  +0.01%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "309.95"; ForceText = $true }
    @{ Cell = "E5"; Value = "  -0.01%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "104.66"; ForceText = $true }
    @{ Cell = "E6"; Value = "  +3.66%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "0.525"; ForceText = $true }
    @{ Cell = "E7"; Value = "  -2.21%  "; ForceText = $false }
    @{ Cell = "E8"; Value = "  +0.05%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "0.520"; ForceText = $true }
    @{ Cell = "E9"; Value = "  +1.83%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "36.23"; ForceText = $true }
    @{ Cell = "E10"; Value = "  +0.32%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "52.82"; ForceText = $true }
    @{ Cell = "E11"; Value = "  +1.08%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "0.0813"; ForceText = $true }
    @{ Cell = "E12"; Value = "  -1.00%  "; ForceText = $false }
    @{ Cell = "E13"; Value = "  -0.80%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "7.01"; ForceText = $true }
    @{ Cell = "E14"; Value = "  +0.22%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "2.742.46"; ForceText = $true }
    @{ Cell = "E15"; Value = "  +3.06%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "15.72"; ForceText = $true }
    @{ Cell = "E16"; Value = "  +5.36%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "2.375.42"; ForceText = $true }
    @{ Cell = "E17"; Value = "  +3.46%  "; ForceText = $false }
    @{ Cell = "E18"; Value = "  +1.70%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "43.347.66"; ForceText = $true }
    @{ Cell = "E19"; Value = "  +0.58%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "12.01"; ForceText = $true }
    @{ Cell = "E20"; Value = "  -4.61%  "; ForceText = $false }
    @{ Cell = "B21"; Value = "Uniswap"; ForceText = $false }
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; ForceText = $false }
    @{ Cell = "D21"; Value = "6.30"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +3.61%  "; ForceText = $false }
    @{ Cell = "B22"; Value = "ShibaInu"; ForceText = $false }
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; ForceText = $false }
    @{ Cell = "D22"; Value = "0.0₃0927"; ForceText = $true }
    @{ Cell = "E22"; Value = "  +0.83%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "68.47"; ForceText = $true }
    @{ Cell = "E23"; Value = "  +0.42%  "; ForceText = $false }
    @{ Cell = "D24"; Value = "242.39"; ForceText = $true }
    @{ Cell = "E24"; Value = "  +0.82%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "2.06"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +1.48%  "; ForceText = $false }
    @{ Cell = "E26"; Value = "  -0.44%  "; ForceText = $false }
    @{ Cell = "E27"; Value = "  +0.25%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "26.10"; ForceText = $true }
    @{ Cell = "E28"; Value = "  +7.98%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "2.30"; ForceText = $true }
    @{ Cell = "E29"; Value = "  +8.20%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "36.94"; ForceText = $true }
    @{ Cell = "E30"; Value = "  -5.55%  "; ForceText = $false }
    @{ Cell = "E31"; Value = "  -0.55%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "161.92"; ForceText = $true }
    @{ Cell = "E32"; Value = "  -2.19%  "; ForceText = $false }
    @{ Cell = "D33"; Value = "5.28"; ForceText = $true }
    @{ Cell = "E33"; Value = "  -1.18%  "; ForceText = $false }
    @{ Cell = "E34"; Value = "  +0.02%  "; ForceText = $false }
    @{ Cell = "D35"; Value = "18.34"; ForceText = $true }
    @{ Cell = "E35"; Value = "  +2.63%  "; ForceText = $false }
    @{ Cell = "B36"; Value = "WEMIXToken"; ForceText = $false }
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; ForceText = $false }
    @{ Cell = "D36"; Value = "2.54"; ForceText = $true }
    @{ Cell = "E36"; Value = "  +6.50%  "; ForceText = $false }
    @{ Cell = "B37"; Value = "LidoDAOToken"; ForceText = $false }
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; ForceText = $false }
    @{ Cell = "D37"; Value = "3.13"; ForceText = $true }
    @{ Cell = "E37"; Value = "  -1.22%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "0.0741"; ForceText = $true }
    @{ Cell = "E38"; Value = "  -0.17%  "; ForceText = $false }
    @{ Cell = "B39"; Value = "ARBITRUM"; ForceText = $false }
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; ForceText = $false }
    @{ Cell = "D39"; Value = "1.95"; ForceText = $true }
    @{ Cell = "E39"; Value = "  +5.39%  "; ForceText = $false }
    @{ Cell = "B40"; Value = "RenderToken"; ForceText = $false }
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; ForceText = $false }
    @{ Cell = "D40"; Value = "4.66"; ForceText = $true }
    @{ Cell = "E40"; Value = "  +11.23%  "; ForceText = $false }
    @{ Cell = "E41"; Value = "  +0.22%  "; ForceText = $false }
    @{ Cell = "E42"; Value = "  -1.26%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "2.43"; ForceText = $true }
    @{ Cell = "E43"; Value = "  +5.07%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "20.35"; ForceText = $true }
    @{ Cell = "E44"; Value = "  +3.19%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "2.007.03"; ForceText = $true }
    @{ Cell = "E45"; Value = "  +1.91%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "3.20"; ForceText = $true }
    @{ Cell = "E46"; Value = "  +5.76%  "; ForceText = $false }
    @{ Cell = "E47"; Value = "  +0.18%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "10.40"; ForceText = $true }
    @{ Cell = "E48"; Value = "  +5.69%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "57.91"; ForceText = $true }
    @{ Cell = "E49"; Value = "  +4.99%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "2.94"; ForceText = $true }
    @{ Cell = "E50"; Value = "  -2.52%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "2.579.73"; ForceText = $true }
    @{ Cell = "E51"; Value = "  +2.01%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
